$wb = $excel.ActiveWorkbook

$wsDesc  = $wb.Worksheets.Item("experiment_description")
$wsSpec  = $wb.Worksheets.Item("experiment_specification")
$wsRunD  = $wb.Worksheets.Item("run_description")
$wsRunS  = $wb.Worksheets.Item("run_specification")

# --- Experiment 17 ("prufa" / vaccine test) -------------------------------

# experiment_description row 25
$wsDesc.Cells.Item(25, 1).Value = 17
$wsDesc.Cells.Item(25, 2).Value = "prufa"
$wsDesc.Cells.Item(25, 3).Value = "base"
$wsDesc.Cells.Item(25, 4).Value = "1;4"
$wsDesc.Cells.Item(25, 5).Value = "none"
$wsDesc.Cells.Item(25, 6).Value = "none"
$wsDesc.Cells.Item(25, 7).Value = "test"
$wsDesc.Cells.Item(25, 8).Value = "prufa"

# experiment_specification rows 80-82 (home / inpatient_ward / intensive_care_unit)
$wsSpec.Cells.Item(80, 1).Value = 17
$wsSpec.Cells.Item(80, 2).Value = "home"
$wsSpec.Cells.Item(80, 3).Value = "length_of_stay_simple_two_weeks"
$wsSpec.Cells.Item(80, 4).Value = "age_simple_vaccinated"
$wsSpec.Cells.Item(80, 5).Value = "age_simple"

$wsSpec.Cells.Item(81, 1).Value = 17
$wsSpec.Cells.Item(81, 2).Value = "inpatient_ward"
$wsSpec.Cells.Item(81, 3).Value = "none"
$wsSpec.Cells.Item(81, 4).Value = "age_simple_vaccinated"
$wsSpec.Cells.Item(81, 5).Value = "none"

$wsSpec.Cells.Item(82, 1).Value = 17
$wsSpec.Cells.Item(82, 2).Value = "intensive_care_unit"
$wsSpec.Cells.Item(82, 3).Value = "none"
$wsSpec.Cells.Item(82, 4).Value = "none"
$wsSpec.Cells.Item(82, 5).Value = "none"

# --- Runs 19 & 20 (referencing experiments 17 & 18) ------------------------

# run_description rows 24-25
$wsRunD.Cells.Item(24, 1).Value = 19
$wsRunD.Cells.Item(24, 2).Value = 17

$wsRunD.Cells.Item(25, 1).Value = 20
$wsRunD.Cells.Item(25, 2).Value = "19, but time_dependent splitting 10 days"

# --- Experiment 18 ("ten days" vaccine test) -------------------------------

# experiment_description row 26
$wsDesc.Cells.Item(26, 1).Value = 18
$wsDesc.Cells.Item(26, 2).Value = "ten days"
$wsDesc.Cells.Item(26, 3).Value = "base"
$wsDesc.Cells.Item(26, 4).Value = "1;4"
$wsDesc.Cells.Item(26, 5).Value = "none"
$wsDesc.Cells.Item(26, 6).Value = "none"
$wsDesc.Cells.Item(26, 7).Value = "10 days"
$wsDesc.Cells.Item(26, 8).Value = "10 daga split"

# experiment_specification rows 83-85 (home / inpatient_ward / intensive_care_unit)
$wsSpec.Cells.Item(83, 1).Value = 18
$wsSpec.Cells.Item(83, 2).Value = "home"
$wsSpec.Cells.Item(83, 3).Value = "length_of_stay_simple_ten_days"
$wsSpec.Cells.Item(83, 4).Value = "age_simple_vaccinated"
$wsSpec.Cells.Item(83, 5).Value = "age_simple"

$wsSpec.Cells.Item(84, 1).Value = 18
$wsSpec.Cells.Item(84, 2).Value = "inpatient_ward"
$wsSpec.Cells.Item(84, 3).Value = "none"
$wsSpec.Cells.Item(84, 4).Value = "age_simple_vaccinated"
$wsSpec.Cells.Item(84, 5).Value = "none"

$wsSpec.Cells.Item(85, 1).Value = 18
$wsSpec.Cells.Item(85, 2).Value = "intensive_care_unit"
$wsSpec.Cells.Item(85, 3).Value = "none"
$wsSpec.Cells.Item(85, 4).Value = "none"
$wsSpec.Cells.Item(85, 5).Value = "none"

# run_specification rows 49-50
$wsRunS.Cells.Item(49, 1).Value = 19
$wsRunS.Cells.Item(49, 2).Value = 17

$wsRunS.Cells.Item(50, 1).Value = 20
$wsRunS.Cells.Item(50, 2).Value = 18

# --- Selection / active sheet bookkeeping ----------------------------------

$wsDesc.Range("H26").Select() | Out-Null
$wsSpec.Range("C83").Select() | Out-Null
$wsSpec.Activate()
